$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Handback is complete: status text changes everywhere it appears ---
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value     = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value     = "Handed back: in sync with en-US"

# --- New handback timestamps recorded for each locale ---
$wsZhCn.Range("K2").Value = "2016-08-28 00:46:55"
$wsDeDe.Range("K2").Value = "2016-08-28 00:47:02"

# --- Handback is now current, so the stale-version error detail clears ---
$wsZhCn.Range("P2").Value = ""
$wsDeDe.Range("P2").Value = ""

# --- Column widths follow the wider/narrower new cell content (report regenerated) ---
$wsOverview.Columns.Item(5).ColumnWidth  = 29.166666666666664
$wsOverview.Columns.Item(6).ColumnWidth  = 29.166666666666664
$wsZhCn.Columns.Item(3).ColumnWidth      = 29.166666666666664
$wsZhCn.Columns.Item(16).ColumnWidth     = 12.833333333333332
$wsDeDe.Columns.Item(3).ColumnWidth      = 29.166666666666664
$wsDeDe.Columns.Item(16).ColumnWidth     = 12.833333333333332
